$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; E=0.9926147277508489; L=0.9955398523336033},
    @{Row=3; E=0.9936372048519304; L=0.9963617723202692},
    @{Row=4; E=0.9942998659930995; L=0.9968940712668345},
    @{Row=5; E=0.9945786998346017; L=0.997117960005301},
    @{Row=6; E=0.9946255319796338; L=0.9971555583673453},
    @{Row=7; E=0.9943035907982488; L=0.9968970624462087},
    @{Row=8; E=0.9929600610674301; L=0.995817528259106},
    @{Row=9; E=0.9906006454969559; L=0.9939188001724441},
    @{Row=10; E=0.989033133672735; L=0.9926553831429383},
    @{Row=11; E=0.988355674866747; L=0.9921088820399291},
    @{Row=12; E=0.9881042295826724; L=0.9919059725120875},
    @{Row=13; E=0.9881581567098651; L=0.9919494934313052},
    @{Row=14; E=0.9883348863814464; L=0.9920921077337197},
    @{Row=15; E=0.9884438009545853; L=0.9921799884222134},
    @{Row=16; E=0.9890781214508737; L=0.9926916645766087},
    @{Row=17; E=0.989476357848556; L=0.9930127773699352},
    @{Row=18; E=0.9897087662937556; L=0.9932001317071769},
    @{Row=19; E=0.9897880325774034; L=0.9932640239640975},
    @{Row=20; E=0.9894336180360679; L=0.9929783193494216},
    @{Row=21; E=0.9882828385668253; L=0.9920501090198105},
    @{Row=22; E=0.9875604150241495; L=0.9914670000341481},
    @{Row=23; E=0.9879432794643023; L=0.991776070289318},
    @{Row=24; E=0.9894529299347244; L=0.9929938892766442},
    @{Row=25; E=0.9912096547607049; L=0.9944092447426414}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 12).Value = $u.L
}

